$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 521
$ws1.Range("F7").Value = 179
$ws1.Range("F9").Value = 990
$ws1.Range("F11").Value = 221
$ws1.Range("F15").Value = 266
$ws1.Range("F16").Value = 571
$ws1.Range("F18").Value = 1315
$ws1.Range("F20").Value = 443
$ws1.Range("F21").Value = 1143
$ws1.Range("F22").Value = 2831
$ws1.Range("F23").Value = 1364
$ws1.Range("F26").Value = 1259
$ws1.Range("F29").Value = 338
$ws1.Range("F30").Value = 2279
$ws1.Range("F31").Value = 337
$ws1.Range("F32").Value = 297
$ws1.Range("F33").Value = 1370

# Sheet "演出" (sheet2): update column F value
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 358

# Sheet "全部类型" (sheet4): update column F values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 521
$ws4.Range("F9").Value = 358
$ws4.Range("F13").Value = 179
$ws4.Range("F16").Value = 990
$ws4.Range("F18").Value = 221
$ws4.Range("F27").Value = 266
$ws4.Range("F28").Value = 571
$ws4.Range("F30").Value = 1315
$ws4.Range("F32").Value = 443
$ws4.Range("F33").Value = 1143
$ws4.Range("F34").Value = 2831
$ws4.Range("F35").Value = 1364
$ws4.Range("F38").Value = 1259
$ws4.Range("F43").Value = 338
$ws4.Range("F44").Value = 2279
$ws4.Range("F45").Value = 337
$ws4.Range("F46").Value = 297
$ws4.Range("F47").Value = 1370
